$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 63 data (date 2019-04-23 / Excel serial 43578)
$rowNum = 63
$values = @(43578, 1, 1, 5, 22, 3, 1, 1, 512, 1, 1, 1, 2, 1, 2, 1, 0, 2, 2, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($rowNum, $col).Value = $values[$i]
}

# Column A is formatted as a date (style matches existing date column above)
$ws.Cells.Item($rowNum - 1, 1).Copy()
$ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($rowNum, 1).Value = $values[0]

# Update selection to match the target state
$ws.Range("I64").Select()
